# Actualización automática del mapa (2025-10-13 07:27:33)
# Appends 4 new incident rows (46-49) to the single "Optical_Power" sheet,
# mirroring the columns: Caso, F. De Reclamo, Direccion, Comuna, OT,
# Proveedor Asignado, Estado, Observaciones, Attachments, Tipo de tarea,
# Equipo, Tipo de Elemento, Coordenada_X, Coordenada_Y, Operacion, Zona, PD, N2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based) that must stay numeric: I=9 (Attachments), M=13 (Coordenada_X), N=14 (Coordenada_Y)
$numericCols = @(9, 13, 14)

$rows = @(
    @("6475", "9/17/2025", "Av Amancio Alcorta 3570", "4", "809800213", "Optical Power", "Pendiente", "aplomar", 1, "Aplomo", "Sin equipos", "Terminal", -58.409278, -34.653566, "San Telmo", "Capital Sur", "PPT-H", "Fuera de Poligono OVL"),
    @("-602", "9/18/2025", "Agustin de vedia 2110", "7", "809837501", "Optical Power", "Pendiente", "Picada", 1, "Cambio", "Sin equipos", "Pasante", -58.435679, -34.643992, "Boedo", "Capital Sur", "PPT-O", "Fuera de Poligono OVL"),
    @("7296", "9/24/2025", "VEDIA, AGUSTIN DE 2130", "7", "809979719", "Optical Power", "Pendiente", "3 columnas picadas ver en calle cuales son necesarias cambiar con Pablo ", 1, "Cambio", "Sin equipos", "Pasante", -58.435634, -34.64412, "Boedo", "Capital Sur", "PPT-O", "Fuera de Poligono OVL"),
    @("-634", "10/8/2025", "Curapaligue 1127", "7", "ICD31464856", "Optical Power", "Pendiente", "Colocar columna donde se marca en la foto pasante 150 o 200", 1, "Cambio", "Sin equipos", "Pasante", -58.446624, -34.635851, "Boedo", "Capital Sur", "PPT-M", "Fuera de Poligono OVL")
)

$startRow = 46
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($numericCols -contains $c) {
            # Numeric columns: write as real numbers.
            $cell.Value = $data[$c - 1]
        } else {
            # Force text storage (some values, e.g. "6475", "9/17/2025", "4",
            # "809800213", look numeric/date-like and would otherwise be
            # auto-converted). Format as Text, assign, then drop back to the
            # Normal style so no explicit formatting lingers on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $data[$c - 1]
            $cell.Style = "Normal"
        }
    }
}
